# Implements: "implemented choice of scaling - minmax vs standard"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Backfill File column for two already-existing rows ---
$ws.Range("A40").Value = "2023-03-10-1727_RF_scaled_avg.csv"
$ws.Range("A41").Value = "2023-03-10-1724_LogReg_scaled_avg.csv"

# --- Row 44: RandomForest, minmax scaling center data before avg ---
$ws.Range("B44").Value = "RandomForest"
$ws.Range("C44").Value = "MoCo"
$ws.Range("D44").Value = "Centers"
$ws.Range("E44").Value = "1 x 3"
$ws.Range("F44").Value = "average"
$ws.Range("I44").Value = "min max scaling center data before avg"
$ws.Range("J44").Value = 0.604
$ws.Range("K44").Value = 0.744
$ws.Range("L44").Value = 0.616
$ws.Range("M44").Value = "0.654 (0.064)"

# --- Row 45: RandomForest, standard scaling center data before avg ---
$ws.Range("B45").Value = "RandomForest"
$ws.Range("C45").Value = "MoCo"
$ws.Range("D45").Value = "Centers"
$ws.Range("E45").Value = "1 x 3"
$ws.Range("F45").Value = "average"
$ws.Range("I45").Value = "standard scaling center data before avg"
$ws.Range("J45").Value = 0.611
$ws.Range("K45").Value = 0.625
$ws.Range("L45").Value = 0.648
$ws.Range("M45").Value = "0.628 (0.015)"

# --- Row 46: minmax submission (RandomForest) ---
$ws.Range("A46").Value = "2023-03-10-1836_RF_minmax_avg.csv"
$ws.Range("B46").Value = "RandomForest"
$ws.Range("C46").Value = "MoCo"
$ws.Range("D46").Value = "Centers"
$ws.Range("E46").Value = "1 x 3"
$ws.Range("F46").Value = "average"
$ws.Range("I46").Value = "min max scaling center data before avg"

# --- Row 47: minmax submission (Logistic Classifier) ---
$ws.Range("A47").Value = "2023-03-10-1838_LogReg_minmax_avg.csv"
$ws.Range("B47").Value = "Logistic Classifier"
$ws.Range("C47").Value = "MoCo"
$ws.Range("D47").Value = "Centers"
$ws.Range("E47").Value = "1 x 3"
$ws.Range("F47").Value = "average"
$ws.Range("I47").Value = "min max scaling center data before avg"

# --- Grow the table to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:O47"))

# --- Extend the "top 10 / top 2" conditional formatting ranges that tracked the table body ---
$rK = $ws.Range("K2:K43")
for ($i = 1; $i -le $rK.FormatConditions.Count; $i++) {
    $fc = $rK.FormatConditions.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$K`$2:`$K`$43") {
        $fc.ModifyAppliesToRange($ws.Range("K2:K47"))
    }
}
$rJ = $ws.Range("J2:J43")
for ($i = 1; $i -le $rJ.FormatConditions.Count; $i++) {
    $fc = $rJ.FormatConditions.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$J`$2:`$J`$43") {
        $fc.ModifyAppliesToRange($ws.Range("J2:J47"))
    }
}

# --- Update viewport / selection to match where the new data was entered ---
$ws.Range("J46").Select()
try {
    $excel.ActiveWindow.ScrollRow = 22
    $excel.ActiveWindow.ScrollColumn = 2
} catch {
    # viewport scroll position is cosmetic only; ignore if unsupported
}
